$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: offset id values by 52500 (1..10 -> 52501..52510)
for ($r = 1; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = 52500 + $r
}

# Fix the misspelled name in row 3
$ws.Range("B3").Value = "Preetika Shetty"

# Column C score corrections
$ws.Range("C1").Value = 80
$ws.Range("C2").Value = 80
$ws.Range("C3").Value = 97
$ws.Range("C4").Value = 96
$ws.Range("C6").Value = 95
$ws.Range("C9").Value = 91
$ws.Range("C10").Value = 90

# Update selection to active cell B3
$ws.Range("B3").Select()
